$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.839.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.64%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '3.689.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +8.26%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''589.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''180.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '3.680.15'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.19%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.622'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.07%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.615'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.81%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''50.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.49%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.0000288'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '4.285.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.27%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''686.45'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''9.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.45%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '3.692.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.48%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '71.887.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''11.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.24%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.946'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.76%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''6.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +17.55%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''17.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.55%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '  +3.42%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '  +3.84%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''2.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.57%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '  +4.84%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''35.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.95%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '  +5.60%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''7.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.84%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''4.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.96%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''570.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '  +2.72%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '  +4.05%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''59.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.67%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '3.796.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.12%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '  +5.33%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '  +5.14%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''35.59'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '  +6.00%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '  +9.72%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '  +4.17%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '  +5.07%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '  +8.56%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''3.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = '''1.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.51%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''134.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.58%  '
$ws.Range("E51").Style = "Normal"
